$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: clone formatting from row 6 (keeps cell styles identical/deduplicated), then fill values ---
$ws.Range("A6:H6").Copy($ws.Range("A7:H7"))

$ws.Range("A7").Value2 = 'Medium'
$ws.Range("B7").Value2 = 'zigzag conversion'
$ws.Range("D7").Value2 = 'NY'
$ws.Range("E7").Value2 = 'link'
$ws.Range("F7").Value2 = 'move pointer, and change direction if necessary '
$ws.Range("G7").Value2 = 'Instead of (I,j) in  a 2D list, sometimes  it is enough to just use one cursor of rows. // Also, try to consider more test cases (e.g. you can randomly generate `s` of different length, and consider `numRows` from 0 to len(s) )'
$ws.Range("H7").Value2 = ""

# --- Row 6 & 7: update/set file names ---
$ws.Range("C6").Value2 = 'longest_sub_palindrom.py'
$ws.Range("C7").Value2 = 'zigzag_convertor.py'

# Hyperlink for E7 (adds relationship + forces a style change on the cell, so we restore
# the cell's formatting/value right after from E6, which keeps styles deduplicated).
$ws.Hyperlinks.Add($ws.Range("E7"), "https://leetcode.com/problems/zigzag-conversion/") | Out-Null
$ws.Range("E6").Copy()
$ws.Range("E7").PasteSpecial(-4122) | Out-Null
$ws.Range("E7").Value2 = 'link'

# --- Row heights ---
$ws.Rows.Item(6).RowHeight = 43.2
$ws.Rows.Item(7).RowHeight = 72

# --- Column widths (col C and col G got wider to fit the new content) ---
$ws.Columns.Item(3).ColumnWidth = 24.333333333333332
$ws.Columns.Item(7).ColumnWidth = 43.166666666666664

# --- View: scroll down one row so row 2 is the top visible row ---
$excel.ActiveWindow.ScrollRow = 2

$excel.CutCopyMode = 0
